# Update with Correct Forecast output
# - Rename Sheet1 -> "Sales vs PO", insert an "Order Week" column holding the
#   original week-ending dates, shift the weekly PO quantity column right and
#   zero it out, and bump the "ds"/"y" dates forward by 6 days.
# - Add three new sheets: "Weekly Growth", "Volume Insights", "Prediction Info"
#   summarizing the PO forecast.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Sales vs PO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Insert a new column C ("Order Week"); this shifts the existing
# PO_Requested_Qty column from C to D (keeping its header style).
$ws1.Columns.Item(3).Insert()
$ws1.Range("C1").Value = "Order Week"

$oldWeekEnding = @(45551,45558,45565,45572,45579,45586,45593,45600,45607,45614,45621,45628,45635,45642,45649)
$newWeekEnding = @(45557,45564,45571,45578,45585,45592,45599,45606,45613,45620,45627,45634,45641,45648,45655)

for ($i = 0; $i -lt $oldWeekEnding.Length; $i++) {
    $row = $i + 2
    $ws1.Range("A$row").Value = $newWeekEnding[$i]
    $ws1.Range("C$row").Value = $oldWeekEnding[$i]
    $ws1.Range("D$row").Value = 0
}

# Give the new "Order Week" data cells the same date format as column A.
$ws1.Range("A2").Copy()
$ws1.Range("C2:C16").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet 2: "Weekly Growth" -- copy Sheet1 (inherits sheetPr/margins/styles),
# then clear its contents and populate fresh data.
# ---------------------------------------------------------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Weekly Growth"
$ws2.Cells.Clear()

$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"

$ws2.Range("A2").Value = 45558
$ws2.Range("B2").Value = 1616
$ws2.Range("C2").Value = 0

$ws2.Range("A3").Value = 45628
$ws2.Range("B3").Value = 144
$ws2.Range("C3").Value = -91.0891089108911

# Match header style (bold + border) and date style used on Sheet 1.
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet 3: "Volume Insights"
# ---------------------------------------------------------------------------
$ws1.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Volume Insights"
$ws3.Cells.Clear()

$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws3.Range("A2").Value = 1760
$ws3.Range("B2").Value = 880
$ws3.Range("C2").Value = 1616
$ws3.Range("D2").Value = 144

$ws1.Range("A1:C1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet 4: "Prediction Info"
# ---------------------------------------------------------------------------
$ws1.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "Prediction Info"
$ws4.Cells.Clear()

$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Range("A2").Value = 0

$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

# Leave the first sheet as the active / selected tab, matching the target.
$ws1.Activate()
